$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the three "Phân công ..." bullet items into a single sentence.
#    "Phân công người thiết kế logo cho nhóm." becomes
#    "Phân công người thiết kế logo cho nhóm, người tạo bản khảo sát,
#     người viết báo cáo."
#    and the two paragraphs that used to hold the other two bullets are
#    removed entirely.
# ---------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Phân công người thiết kế logo cho nhóm")) {
        $target = $cand
        break
    }
}

$tr = $target.Range
# Drop the trailing "." - it sits right before the paragraph mark.
$periodRange = $d.Range($tr.End - 2, $tr.End - 1)
$periodRange.Delete()

# Re-fetch the (now shifted) paragraph range and insert the continuation
# text right before the paragraph mark.
$tr = $target.Range
$insPos = $tr.End - 1
$insertPoint = $d.Range($insPos, $insPos)
$insertPoint.InsertAfter(", người tạo bản khảo sát, người viết báo cáo.")

# The insertion above lands in the same run as the preceding text because
# the run formatting is identical; split it into its own <w:r> by
# round-tripping a character formatting property over just the new text.
$tr = $target.Range
$newEndPos = $tr.End - 1
$newRunRange = $d.Range($insPos, $newEndPos)
$newRunRange.Font.Bold = 1
$newRunRange.Font.Bold = 0

# Now remove the two now-redundant bullet paragraphs
# ("Phân công người tạo bản khảo sát." and "Phân công người viết báo cáo.").
$survey = $null
$report = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t.StartsWith("Phân công người tạo bản khảo sát")) { $survey = $cand }
    if ($t.StartsWith("Phân công người viết báo cáo")) { $report = $cand }
}

$delRange = $d.Range($survey.Range.Start, $report.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the "... làm bản báo
#    cáo." paragraph to the very start of the "Lên kế hoạch khảo sát."
#    paragraph.
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$planParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Lên kế hoạch khảo sát")) {
        $planParagraph = $cand
        break
    }
}

$bmPos = $planParagraph.Range.Start
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
